$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old demo content (rows 1-9 / cols A-C) so the sheet starts clean
$ws.Cells.Clear()

# Write the new full-width editable table row into A2:O2.
# Leading "'" keeps numeric-looking values as text (matching the source data).
$values = @(
    "2025-01-30 16:03:19",
    "'2001",
    "Rental Income",
    "'1100",
    "'2034",
    "'2345",
    "'456",
    "'345",
    "'657",
    "'657",
    "'777",
    "'787",
    "'788",
    "'790",
    "'800"
)

$data = New-Object 'object[,]' 1,15
for ($i = 0; $i -lt $values.Length; $i++) {
    $data[0, $i] = $values[$i]
}

$ws.Range("A2:O2").Value = $data

# Drop the quote-prefix formatting artifact left by the "'" text-entry trick
# so the cells end up with plain default styling, same as the source file.
$ws.Range("A2:O2").ClearFormats()

# Match the recorded selection / active cell from the saved view state
$ws.Range("H18").Select() | Out-Null
